$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 99, shifting existing data (rows 99-129) down to (100-130)
$ws.Rows.Item(99).Insert()

# Fill in the new row 99 with its data
$ws.Cells.Item(99, 1).Value = 10
$ws.Cells.Item(99, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(99, 3).Value = "La Araucanía"
$ws.Cells.Item(99, 4).Value = 44627
$ws.Cells.Item(99, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(99, 5).Value = 9
$ws.Cells.Item(99, 6).Value = "Fruta"
$ws.Cells.Item(99, 7).Value = 100104
$ws.Cells.Item(99, 8).Value = "Frutos de pepita"
$ws.Cells.Item(99, 9).Value = 100104003
$ws.Cells.Item(99, 10).Value = "Membrillo"
$ws.Cells.Item(99, 11).Value = "Champion"
$ws.Cells.Item(99, 12).Value = "Primera"
$ws.Cells.Item(99, 13).Value = 80
$ws.Cells.Item(99, 14).Value = 16000
$ws.Cells.Item(99, 15).Value = 16000
$ws.Cells.Item(99, 16).Value = 16000
$ws.Cells.Item(99, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(99, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(99, 19).Value = 889
$ws.Cells.Item(99, 20).Value = 18
